$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 9, pushing existing rows 9..92 down to 10..93
$ws.Rows.Item(9).Insert()

# Populate the newly inserted row 9 with its data
$ws.Cells.Item(9, 1).Value = 4
$ws.Cells.Item(9, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(9, 3).Value = "Los Lagos"
$ws.Cells.Item(9, 4).Value = 44530
$ws.Cells.Item(9, 5).Value = 10
$ws.Cells.Item(9, 6).Value = 100112022
$ws.Cells.Item(9, 7).Value = "Arveja Verde"
$ws.Cells.Item(9, 8).Value = "Sin especificar"
$ws.Cells.Item(9, 9).Value = "Primera"
$ws.Cells.Item(9, 10).Value = 120
$ws.Cells.Item(9, 11).Value = 18000
$ws.Cells.Item(9, 12).Value = 18000
$ws.Cells.Item(9, 13).Value = 18000
$ws.Cells.Item(9, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(9, 15).Value = "Región del Maule"
$ws.Cells.Item(9, 16).Value = 720
$ws.Cells.Item(9, 17).Value = 25
$ws.Cells.Item(9, 18).Value = "Hortaliza"

# Match the date-number-format style already used by the other rows in column D
$ws.Cells.Item(9, 4).NumberFormat = $ws.Cells.Item(10, 4).NumberFormat
